# The deck currently uses the "Integral" (Red Violet) design theme on its
# single slide master / notes master.  The authored change swaps the theme
# applied to the deck back to the plain built-in "Office Theme" colour
# scheme (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink), i.e. the deck's theme
# part ends up holding the Office Theme colours instead of the Red Violet
# ones.
#
# Helper: pack R/G/B (0-255) into the BGR-packed Long that
# ThemeColor.RGB / RGB() uses throughout the PowerPoint object model.
function ColorRGB($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme

# Office Theme colour scheme values, in the fixed 1-12 ThemeColorScheme
# order: dk1, lt1, dk2, lt2, accent1, accent2, accent3, accent4, accent5,
# accent6, hlink, folHlink.
$officeThemeColors = @(
    (ColorRGB 0x00 0x00 0x00),  # dk1      000000
    (ColorRGB 0xFF 0xFF 0xFF),  # lt1      FFFFFF
    (ColorRGB 0x44 0x54 0x6A),  # dk2      44546A
    (ColorRGB 0xE7 0xE6 0xE6),  # lt2      E7E6E6
    (ColorRGB 0x5B 0x9B 0xD5),  # accent1  5B9BD5
    (ColorRGB 0xED 0x7D 0x31),  # accent2  ED7D31
    (ColorRGB 0xA5 0xA5 0xA5),  # accent3  A5A5A5
    (ColorRGB 0xFF 0xC0 0x00),  # accent4  FFC000
    (ColorRGB 0x44 0x72 0xC4),  # accent5  4472C4
    (ColorRGB 0x70 0xAD 0x47),  # accent6  70AD47
    (ColorRGB 0x05 0x63 0xC1),  # hlink    0563C1
    (ColorRGB 0x95 0x4F 0x72)   # folHlink 954F72
)

for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $colorScheme.Item($i).RGB = $officeThemeColors[$i - 1]
}
